$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)
$nl = [char]11

$tbl.Cell(1,1).Range.Text = "24 x 86" + $nl + "  8    6" + $nl + "  ----" + $nl + "2|    |" + $nl + "4|    |"
$tbl.Cell(1,2).Range.Text = "43 x 33" + $nl + "  3    3" + $nl + "  ----" + $nl + "4|    |" + $nl + "3|    |"
$tbl.Cell(1,3).Range.Text = "12 x 46" + $nl + "  4    6" + $nl + "  ----" + $nl + "1|    |" + $nl + "2|    |"
$tbl.Cell(2,1).Range.Text = "46 x 90" + $nl + "  9    0" + $nl + "  ----" + $nl + "4|    |" + $nl + "6|    |"
$tbl.Cell(2,2).Range.Text = "58 x 71" + $nl + "  7    1" + $nl + "  ----" + $nl + "5|    |" + $nl + "8|    |"
$tbl.Cell(2,3).Range.Text = "84 x 51" + $nl + "  5    1" + $nl + "  ----" + $nl + "8|    |" + $nl + "4|    |"
$tbl.Cell(3,1).Range.Text = "52 x 96" + $nl + "  9    6" + $nl + "  ----" + $nl + "5|    |" + $nl + "2|    |"
$tbl.Cell(3,2).Range.Text = "66 x 36" + $nl + "  3    6" + $nl + "  ----" + $nl + "6|    |" + $nl + "6|    |"
$tbl.Cell(3,3).Range.Text = "53 x 18" + $nl + "  1    8" + $nl + "  ----" + $nl + "5|    |" + $nl + "3|    |"
$tbl.Cell(4,1).Range.Text = "69 x 57" + $nl + "  5    7" + $nl + "  ----" + $nl + "6|    |" + $nl + "9|    |"
$tbl.Cell(4,2).Range.Text = "20 x 65" + $nl + "  6    5" + $nl + "  ----" + $nl + "2|    |" + $nl + "0|    |"
$tbl.Cell(4,3).Range.Text = "57 x 92" + $nl + "  9    2" + $nl + "  ----" + $nl + "5|    |" + $nl + "7|    |"
$tbl.Cell(5,1).Range.Text = "20 x 93" + $nl + "  9    3" + $nl + "  ----" + $nl + "2|    |" + $nl + "0|    |"
$tbl.Cell(5,2).Range.Text = "61 x 94" + $nl + "  9    4" + $nl + "  ----" + $nl + "6|    |" + $nl + "1|    |"
$tbl.Cell(5,3).Range.Text = "42 x 97" + $nl + "  9    7" + $nl + "  ----" + $nl + "4|    |" + $nl + "2|    |"
